$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the title text in A4: "Commercial revenue (2023 USD)" -> "Commercial revenue (millions of 2023 USD)"
$ws.Range("A4").Value = "Commercial revenue (millions of 2023 USD)"

# Update the selected/active cell on the sheet from E8 to A4
$ws.Range("A4").Select()
